$p = $ppt.ActivePresentation

# --- Slide 2: merge the "Baseline " + "Model" runs (inside nested group) into
#     a single run reading "Baseline Model", keeping the first run's formatting
#     and the shape's original size (the autosize box would otherwise shrink
#     its height when the text is rewritten, so we restore it explicitly). ---
$s2 = $p.Slides.Item(2)
$grp2 = $s2.Shapes.Item(3)          # 组合 2 (outer group)
$tb7 = $grp2.GroupItems.Item(1)     # 文本框 7 (flattened through 组合 3)
$tr7 = $tb7.TextFrame.TextRange
$run1 = $tr7.Runs(1)
$run2 = $tr7.Runs(2)
$run1.Text = "Baseline Model"
$run2.Text = ""
$tb7.Height = 33.928031

# --- Slide 3: merge the "Baseline " + "Model" runs of the 2nd paragraph into
#     a single run reading "Baseline Model", keeping the first run's formatting. ---
$s3 = $p.Slides.Item(3)
$tb13 = $s3.Shapes.Item(3)          # 文本框 13
$tr13 = $tb13.TextFrame.TextRange
$para2 = $tr13.Paragraphs(2)
$para2.Text = "placeholder_tmp"
$para2b = $tr13.Paragraphs(2)
$para2b.Text = "Baseline Model"

# --- Slide 4: resize "TextBox 2" shape ---
$s4 = $p.Slides.Item(4)
$tbox2 = $s4.Shapes.Item(6)         # TextBox 2
$tbox2.Left = 48.916064
$tbox2.Top = 173.129769
$tbox2.Width = 601.282441
$tbox2.Height = 72.703152
